$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 3): rename/re-order measurement headers ---
$ws.Range("B3").Value = "Ushunt [mV]"
$ws.Range("C3").Value = "Rshunt [Ohm]"
$ws.Range("D3").Value = "Imote[mA]"
$ws.Range("E3").Value = "Ubattery [V]"
$ws.Range("F3").Value = "Pmote [W]"
$ws.Range("G3").Value = "Time [s]"
$ws.Range("H3").Value = "Energy [J]"

# --- Row 4 (Receiver Idle): add /10 factor to power formula ---
$ws.Range("F4").Formula = "=(E4-(B4*0.001))*(B4*0.001)/10"

# --- Rows 5-11: shared power formula with /10 factor ---
$ws.Range("F5:F11").Formula = "=(E5-(B5*0.001))*(B5*0.001)/10"

# --- Row 8 (Sender Radio turned on): corrected voltage reading ---
$ws.Range("B8").Value = 188.75

# --- Row 10 (Sender sending, no comp): corrected voltage reading ---
$ws.Range("B10").Value = 188.75
$ws.Range("F10").Formula = "=(E10-(B10*0.001))*(B10*0.001)/10"
$ws.Range("H10").Formula = "=F10*G10"

# --- Row 11 (Sender sending, compression): energy formula recomputed ---
$ws.Range("H11").Formula = "=F11*G11"

# --- Rows 12-14: power formula reverts to the original (no /10) ---
$ws.Range("F12:F14").Formula = "=(E12-(B12*0.001))*(B12*0.001)"

# --- Selection moves to B10 ---
$ws.Range("B10").Select()
